# Update the cryptos list worksheet with refreshed price/volume data
# (mirrors the "Updated cryptos list" GitHub Actions commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.353.41'
$ws.Range("E2").Value = '  -3.11%  '
$ws.Range("D3").Value = '3.554.28'
$ws.Range("E3").Value = '  -3.55%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''585.33'
$ws.Range("E5").Value = '  -1.44%  '
$ws.Range("D6").Value = '''179.89'
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("D7").Value = '3.546.72'
$ws.Range("E7").Value = '  -3.62%  '
$ws.Range("D8").Value = '''0.605'
$ws.Range("E8").Value = '  -4.07%  '
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").Value = '''0.666'
$ws.Range("E10").Value = '  -6.92%  '
$ws.Range("D11").Value = '''0.143'
$ws.Range("E11").Value = '  -11.36%  '
$ws.Range("D12").Value = '''53.19'
$ws.Range("E12").Value = '  -4.89%  '
$ws.Range("D13").Value = '''0.0000249'
$ws.Range("E13").Value = '  -14.84%  '
$ws.Range("D14").Value = '''9.73'
$ws.Range("E14").Value = '  -8.99%  '
$ws.Range("D15").Value = '4.121.17'
$ws.Range("E15").Value = '  -3.32%  '
$ws.Range("D16").Value = '3.542.65'
$ws.Range("E16").Value = '  -3.66%  '
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").Value = '''18.23'
$ws.Range("E18").Value = '  -5.95%  '
$ws.Range("D19").Value = '66.012.10'
$ws.Range("E19").Value = '  -3.27%  '
$ws.Range("D20").Value = '''12.01'
$ws.Range("E20").Value = '  -6.50%  '
$ws.Range("E21").Value = '  -7.16%  '
$ws.Range("D22").Value = '''390.52'
$ws.Range("E22").Value = '  -4.96%  '
$ws.Range("D23").Value = '''4.26'
$ws.Range("E23").Value = '  -7.22%  '
$ws.Range("D24").Value = '''84.09'
$ws.Range("E24").Value = '  -5.23%  '
$ws.Range("D25").Value = '''2.85'
$ws.Range("E25").Value = '  -5.96%  '
$ws.Range("D26").Value = '''12.15'
$ws.Range("E26").Value = '  -4.19%  '
$ws.Range("D27").Value = '''6.01'
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").Value = '''10.18'
$ws.Range("E28").Value = '  -5.48%  '
$ws.Range("D29").Value = '''3.55'
$ws.Range("E29").Value = '  -8.72%  '
$ws.Range("D30").Value = '''8.83'
$ws.Range("E30").Value = '  -8.02%  '
$ws.Range("D31").Value = '''30.89'
$ws.Range("E31").Value = '  -5.98%  '
$ws.Range("D32").Value = '''6.68'
$ws.Range("E32").Value = '  -9.53%  '
$ws.Range("D33").Value = '''64.98'
$ws.Range("D34").Value = '''11.84'
$ws.Range("E34").Value = '  -4.79%  '
$ws.Range("D35").Value = '''594.29'
$ws.Range("E35").Value = '  -1.69%  '
$ws.Range("D36").Value = '''0.111'
$ws.Range("E36").Value = '  -7.11%  '
$ws.Range("D37").Value = '''41.09'
$ws.Range("E37").Value = '  -5.35%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("D39").Value = '''0.998'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Value = '''0.369'
$ws.Range("E40").Value = '  -8.13%  '
$ws.Range("D41").Value = '0.0₃0730'
$ws.Range("E41").Value = '  -17.51%  '
$ws.Range("D42").Value = '''0.128'
$ws.Range("E42").Value = '  -6.89%  '
$ws.Range("D43").Value = '2.875.33'
$ws.Range("E43").Value = '  +5.67%  '
$ws.Range("D44").Value = '''2.75'
$ws.Range("E44").Value = '  -9.44%  '
$ws.Range("D45").Value = '''0.0404'
$ws.Range("E45").Value = '  -8.66%  '
$ws.Range("D46").Value = '''2.39'
$ws.Range("E46").Value = '  -9.80%  '
$ws.Range("E47").Value = '  -4.32%  '
$ws.Range("E48").Value = '  -2.81%  '
$ws.Range("D49").Value = '''2.50'
$ws.Range("E49").Value = '  -8.38%  '
$ws.Range("D50").Value = '''135.39'
$ws.Range("E50").Value = '  -3.33%  '
$ws.Range("D51").Value = '''8.17'
$ws.Range("E51").Value = '  -9.22%  '
